# Actualiza las lecturas de los sensores con los nuevos datos capturados
# (visualizacion en tiempo real + guardado en SQLite).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cada fila: Fila, Timestamp, Temperatura, Humedad, Presion, Luz
$sensorData = @(
    @(2, "2025-06-08 01:48:38", 27.5, 55, 1010.26, 15.83),
    @(3, "2025-06-08 01:48:40", 27.1, 54, 1010.27, 16.67),
    @(4, "2025-06-08 01:48:42", 27.1, 54, 1010.26, 18.33),
    @(5, "2025-06-08 01:48:45", 27.1, 54, 1010.28, 19.17),
    @(6, "2025-06-08 01:48:47", 27.1, 54, 1010.26, 19.17),
    @(7, "2025-06-08 01:48:49", 27.1, 54, 1010.24, 19.17),
    @(8, "2025-06-08 01:48:51", 27.1, 54, 1010.25, 19.17),
    @(9, "2025-06-08 01:48:53", 27.1, 55, 1010.27, 19.17),
    @(10, "2025-06-08 01:48:55", 27.1, 55, 1010.21, 19.17),
    @(11, "2025-06-08 01:48:57", 27.1, 55, 1010.25, 17.5),
    @(12, "2025-06-08 01:48:59", 27.1, 55, 1010.23, 17.5),
    @(13, "2025-06-08 01:49:01", 27.1, 55, 1010.24, 18.33),
    @(14, "2025-06-08 01:49:03", 26.7, 55, 1010.24, 16.67),
    @(15, "2025-06-08 01:49:05", 26.7, 55, 1010.25, 17.5),
    @(16, "2025-06-08 01:49:07", 26.7, 55, 1010.28, 17.5),
    @(17, "2025-06-08 01:49:09", 26.7, 55, 1010.29, 17.5),
    @(18, "2025-06-08 01:49:11", 26.7, 55, 1010.3, 17.5),
    @(19, "2025-06-08 01:49:13", 26.7, 55, 1010.32, 350.83),
    @(20, "2025-06-08 01:49:15", 26.7, 55, 1010.2, 17.5),
    @(21, "2025-06-08 01:49:17", 26.7, 55, 1010.3, 18.33)
)

foreach ($row in $sensorData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
